$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "V" indicator column inserted into column A for every data row (2-41)
$ws.Range("A2:A41").Value = "V"

# Columns F (6) and H (8) are now hidden (col G stays visible, split out from the old merged 7-8 width group)
$ws.Columns.Item(6).Hidden = $true
$ws.Columns.Item(8).Hidden = $true

# Turn on AutoFilter over the full table range
$ws.Range("A1:L41").AutoFilter() | Out-Null

# AutoFilter registers a hidden sheet-scoped defined name (_FilterDatabase)
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='材料資料片 TR2 VHR'!`$A`$1:`$L`$41")
$filterName.Visible = $false

# Selection moved to L9
$ws.Range("L9").Select() | Out-Null
